$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the score input cells (C:G and J:N) for rows 2-6, keeping their styles.
# This mirrors the author's edit where values were deleted but formatting stayed,
# so the dependent SUM formulas (H/O/V/AC/AJ columns) all recompute to 0.
$ws.Range("C2:G6").ClearContents()
$ws.Range("J2:N6").ClearContents()

# Leave selection where the user ended up after clearing the range.
$ws.Range("K13").Select()
